$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 533542.4399999999
$ws.Range("I15").Value = 533542.4399999999
$ws.Range("K15").Value = 1600627.32
$ws.Range("M15").Value = -1600458.32
$ws.Range("H33").Value = 2086.6316
$ws.Range("I33").Value = 380.33334
$ws.Range("J33").Value = 8485.25
$ws.Range("K33").Value = 380.33334
$ws.Range("L33").Value = 8485.25
$ws.Range("M33").Value = -151.33334
$ws.Range("N33").Value = -8943.25
$ws.Range("H55").Value = 624.5
$ws.Range("I55").Value = 450
$ws.Range("K55").Value = 450
$ws.Range("M55").Value = -236
$ws.Range("H62").Value = 7930.375
$ws.Range("I62").Value = 5990.6665
$ws.Range("K62").Value = 5990.6665
$ws.Range("M62").Value = -5366.6665
$ws.Range("H64").Value = 10581.556
$ws.Range("I64").Value = 6664
$ws.Range("K64").Value = 6664
$ws.Range("M64").Value = -6416
$ws.Range("H65").Value = 7930.375
$ws.Range("I65").Value = 5990.6665
$ws.Range("K65").Value = 29953.3325
$ws.Range("M65").Value = -26833.3325
$ws.Range("H67").Value = 10581.556
$ws.Range("I67").Value = 6664
$ws.Range("K67").Value = 6664
$ws.Range("M67").Value = -5806
$ws.Range("H68").Value = 69999
$ws.Range("J68").Value = 69999
$ws.Range("L68").Value = 69999
$ws.Range("N68").Value = -71497
$ws.Range("H71").Value = 69999
$ws.Range("J71").Value = 69999
$ws.Range("L71").Value = 209997
$ws.Range("N71").Value = -217485
$ws.Range("H107").Value = 773.9231
$ws.Range("I107").Value = 587.9167
$ws.Range("K107").Value = 587.9167
$ws.Range("M107").Value = 1332.0833
$ws.Range("H112").Value = 6621.727
$ws.Range("I112").Value = 2444.5
$ws.Range("J112").Value = 7550
$ws.Range("K112").Value = 7333.5
$ws.Range("L112").Value = 22650
$ws.Range("M112").Value = -6225.5
$ws.Range("N112").Value = -24866
$ws.Range("H115").Value = 0
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("M115").ClearContents()
$ws.Range("N115").ClearContents()
$ws.Range("H125").Value = 83446
$ws.Range("J125").Value = 132099.88
$ws.Range("L125").Value = 1188898.92
$ws.Range("N125").Value = -1193818.92
$ws.Range("H131").Value = 3000
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 3000
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 9000
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -19080
$ws.Range("H132").Value = 3398.9216
$ws.Range("I132").Value = 2359.4565
$ws.Range("K132").Value = 7078.369499999999
$ws.Range("M132").Value = -4548.369499999999
$ws.Range("H134").Value = 112499
$ws.Range("J134").Value = 112499
$ws.Range("L134").Value = 112499
$ws.Range("N134").Value = -122639

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 729.85
$ws.Range("J2").Value = 670.1429000000001
$ws.Range("L2").Value = 670.1429000000001
$ws.Range("N2").Value = -896.1429000000001
$ws.Range("H32").Value = 9708.612999999999
$ws.Range("I32").Value = 5461.147
$ws.Range("J32").Value = 24150
$ws.Range("K32").Value = 5461.147
$ws.Range("L32").Value = 24150
$ws.Range("M32").Value = -5174.147
$ws.Range("N32").Value = -24724
$ws.Range("H43").Value = 7331.3335
$ws.Range("J43").Value = 7331.3335
$ws.Range("L43").Value = 7331.3335
$ws.Range("N43").Value = -7957.3335
$ws.Range("H45").Value = 3753.182
$ws.Range("I45").Value = 3445.5
$ws.Range("K45").Value = 3445.5
$ws.Range("M45").Value = -3068.5
$ws.Range("H61").Value = 2828.6584
$ws.Range("I61").Value = 2352.4
$ws.Range("K61").Value = 2352.4
$ws.Range("M61").Value = -2140.4
$ws.Range("H74").Value = 4288.2915
$ws.Range("J74").Value = 7340.769
$ws.Range("L74").Value = 7340.769
$ws.Range("N74").Value = -9088.769
$ws.Range("H77").Value = 4288.2915
$ws.Range("J77").Value = 7340.769
$ws.Range("L77").Value = 36703.845
$ws.Range("N77").Value = -45439.845
$ws.Range("H88").Value = 3995.2727
$ws.Range("I88").Value = 5249.5
$ws.Range("J88").Value = 3278.5715
$ws.Range("K88").Value = 5249.5
$ws.Range("L88").Value = 3278.5715
$ws.Range("M88").Value = -4843.5
$ws.Range("N88").Value = -4090.5715
$ws.Range("H91").Value = 3995.2727
$ws.Range("I91").Value = 5249.5
$ws.Range("J91").Value = 3278.5715
$ws.Range("K91").Value = 5249.5
$ws.Range("L91").Value = 3278.5715
$ws.Range("M91").Value = -3845.5
$ws.Range("N91").Value = -6086.5715
$ws.Range("H110").Value = 989.53845
$ws.Range("I110").Value = 989.53845
$ws.Range("K110").Value = 989.53845
$ws.Range("M110").Value = 1055.46155
$ws.Range("H116").Value = 729.85
$ws.Range("J116").Value = 670.1429000000001
$ws.Range("L116").Value = 670.1429000000001
$ws.Range("N116").Value = -5258.1429
$ws.Range("H122").Value = 3622.1765
$ws.Range("I122").Value = 4263.625
$ws.Range("J122").Value = 3052
$ws.Range("K122").Value = 12790.875
$ws.Range("L122").Value = 9156
$ws.Range("M122").Value = -10340.875
$ws.Range("N122").Value = -14056
$ws.Range("H132").Value = 1712.5555
$ws.Range("I132").Value = 1364.125
$ws.Range("K132").Value = 4092.375
$ws.Range("M132").Value = -1562.375
$ws.Range("H136").Value = 2828.6584
$ws.Range("I136").Value = 2352.4
$ws.Range("K136").Value = 7057.200000000001
$ws.Range("M136").Value = -4507.200000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 729.85
$ws.Range("J3").Value = 670.1429000000001
$ws.Range("L3").Value = 670.1429000000001
$ws.Range("N3").Value = -898.1429000000001
$ws.Range("H5").Value = 23493.2
$ws.Range("I5").Value = 1933.3334
$ws.Range("J5").Value = 55833
$ws.Range("K5").Value = 1933.3334
$ws.Range("L5").Value = 55833
$ws.Range("M5").Value = -1820.3334
$ws.Range("N5").Value = -56059
$ws.Range("H7").Value = 4800
$ws.Range("I7").Value = 3625
$ws.Range("J7").Value = 9500
$ws.Range("K7").Value = 3625
$ws.Range("L7").Value = 9500
$ws.Range("M7").Value = -3512
$ws.Range("N7").Value = -9726
$ws.Range("H26").Value = 14998
$ws.Range("I26").Value = 14998
$ws.Range("K26").Value = 14998
$ws.Range("M26").Value = -14706
$ws.Range("H86").Value = 4679.857
$ws.Range("I86").Value = 1667
$ws.Range("J86").Value = 22757
$ws.Range("K86").Value = 1667
$ws.Range("L86").Value = 22757
$ws.Range("M86").Value = -544
$ws.Range("N86").Value = -25003
$ws.Range("H89").Value = 4679.857
$ws.Range("I89").Value = 1667
$ws.Range("J89").Value = 22757
$ws.Range("K89").Value = 8335
$ws.Range("L89").Value = 113785
$ws.Range("M89").Value = -2719
$ws.Range("N89").Value = -125017
$ws.Range("H107").Value = 1604.9524
$ws.Range("I107").Value = 1078.8462
$ws.Range("K107").Value = 1078.8462
$ws.Range("M107").Value = 841.1538

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1132.1428
$ws.Range("I16").Value = 987.6667
$ws.Range("J16").Value = 1999
$ws.Range("K16").Value = 987.6667
$ws.Range("L16").Value = 1999
$ws.Range("M16").Value = -700.6667
$ws.Range("N16").Value = -2573
$ws.Range("H31").Value = 7504.069
$ws.Range("I31").Value = 4139.5137
$ws.Range("K31").Value = 4139.5137
$ws.Range("M31").Value = -3844.5137
$ws.Range("H34").Value = 7504.069
$ws.Range("I34").Value = 4139.5137
$ws.Range("K34").Value = 4139.5137
$ws.Range("M34").Value = -3937.5137
$ws.Range("H58").Value = 3155.3333
$ws.Range("I58").Value = 2955.5
$ws.Range("K58").Value = 2955.5
$ws.Range("M58").Value = -2752.5
$ws.Range("H94").Value = 1593.9166
$ws.Range("I94").Value = 707.5
$ws.Range("K94").Value = 707.5
$ws.Range("M94").Value = -256.5
$ws.Range("H99").Value = 12251.818
$ws.Range("I99").Value = 15425.125
$ws.Range("K99").Value = 15425.125
$ws.Range("M99").Value = -13927.125
$ws.Range("H103").Value = 30409.2
$ws.Range("J103").Value = 46499.5
$ws.Range("L103").Value = 46499.5
$ws.Range("N103").Value = -48843.5
$ws.Range("H109").Value = 63095
$ws.Range("I109").Value = 45000
$ws.Range("J109").Value = 99285
$ws.Range("K109").Value = 45000
$ws.Range("L109").Value = 99285
$ws.Range("M109").Value = -43960
$ws.Range("N109").Value = -101365
$ws.Range("H113").Value = 1132.1428
$ws.Range("I113").Value = 987.6667
$ws.Range("J113").Value = 1999
$ws.Range("K113").Value = 987.6667
$ws.Range("L113").Value = 1999
$ws.Range("M113").Value = 1182.3333
$ws.Range("N113").Value = -6339
$ws.Range("H126").Value = 12251.818
$ws.Range("I126").Value = 15425.125
$ws.Range("K126").Value = 46275.375
$ws.Range("M126").Value = -43805.375
$ws.Range("H132").Value = 5758.5557
$ws.Range("I132").Value = 3567.4736
$ws.Range("K132").Value = 10702.4208
$ws.Range("M132").Value = -8172.4208
$ws.Range("H134").Value = 2399.721
$ws.Range("I134").Value = 2081.6667
$ws.Range("K134").Value = 6245.000100000001
$ws.Range("M134").Value = -3710.000100000001
$ws.Range("H136").Value = 3155.3333
$ws.Range("I136").Value = 2955.5
$ws.Range("K136").Value = 8866.5
$ws.Range("M136").Value = -6316.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 1772.125
$ws.Range("I69").Value = 863.1667
$ws.Range("K69").Value = 2589.5001
$ws.Range("M69").Value = -1778.5001
$ws.Range("H72").Value = 1772.125
$ws.Range("I72").Value = 863.1667
$ws.Range("K72").Value = 7768.5003
$ws.Range("M72").Value = -3712.5003
$ws.Range("H107").Value = 435.6154
$ws.Range("J107").Value = 570.5
$ws.Range("L107").Value = 1711.5
$ws.Range("N107").Value = -5551.5
$ws.Range("H114").Value = 2459
$ws.Range("I114").Value = 2459
$ws.Range("J114").Value = 0
$ws.Range("K114").Value = 7377
$ws.Range("L114").Value = 0
$ws.Range("M114").Value = -4123
$ws.Range("N114").ClearContents()
$ws.Range("H117").Value = 0
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("M117").ClearContents()
$ws.Range("N117").ClearContents()
$ws.Range("H119").Value = 5677.5
$ws.Range("I119").Value = 4086.2222
$ws.Range("K119").Value = 12258.6666
$ws.Range("M119").Value = -7420.6666
$ws.Range("H120").Value = 13397.308
$ws.Range("I120").Value = 6356.1113
$ws.Range("J120").Value = 29240
$ws.Range("K120").Value = 19068.3339
$ws.Range("L120").Value = 87720
$ws.Range("M120").Value = -14230.3339
$ws.Range("N120").Value = -97396
$ws.Range("H128").Value = 270164.5
$ws.Range("I128").Value = 270164.5
$ws.Range("K128").Value = 810493.5
$ws.Range("M128").Value = -805513.5
$ws.Range("H137").Value = 3228.5881
$ws.Range("I137").Value = 2062.4
$ws.Range("J137").Value = 4894.5713
$ws.Range("K137").Value = 6187.200000000001
$ws.Range("L137").Value = 14683.7139
$ws.Range("M137").Value = -1087.200000000001
$ws.Range("N137").Value = -24883.7139
$ws.Range("H140").Value = 1823.6666
$ws.Range("I140").Value = 1142
$ws.Range("J140").Value = 2018.4286
$ws.Range("K140").Value = 3426
$ws.Range("L140").Value = 6055.2858
$ws.Range("M140").Value = 1754
$ws.Range("N140").Value = -16415.2858

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()
$ws.Range("H48").Value = 38000
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()
$ws.Range("H102").Value = 2476.1924
$ws.Range("I102").Value = 1648.2778
$ws.Range("K102").Value = 1648.2778
$ws.Range("M102").Value = -26.27780000000007
$ws.Range("H122").Value = 4733.5
$ws.Range("I122").Value = 2153.842
$ws.Range("K122").Value = 6461.526
$ws.Range("M122").Value = -4011.526
$ws.Range("H123").Value = 54997.332
$ws.Range("J123").Value = 54997.332
$ws.Range("L123").Value = 54997.332
$ws.Range("N123").Value = -59897.332
$ws.Range("H138").Value = 99997.8
$ws.Range("J138").Value = 99997.8
$ws.Range("L138").Value = 99997.8
$ws.Range("N138").Value = -110277.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 46090.75
$ws.Range("I7").Value = 62617
$ws.Range("J7").Value = 5955.5713
$ws.Range("K7").Value = 62617
$ws.Range("L7").Value = 5955.5713
$ws.Range("M7").Value = -62505
$ws.Range("N7").Value = -6179.5713
$ws.Range("H40").Value = 5183.913
$ws.Range("I40").Value = 3242.5
$ws.Range("J40").Value = 6219.3335
$ws.Range("K40").Value = 3242.5
$ws.Range("L40").Value = 6219.3335
$ws.Range("M40").Value = -3106.5
$ws.Range("N40").Value = -6491.3335
$ws.Range("H46").Value = 3258.5334
$ws.Range("J46").Value = 3567.6924
$ws.Range("L46").Value = 3567.6924
$ws.Range("N46").Value = -3943.6924
$ws.Range("H61").Value = 3803.4443
$ws.Range("I61").Value = 1377
$ws.Range("J61").Value = 10112.2
$ws.Range("K61").Value = 1377
$ws.Range("L61").Value = 10112.2
$ws.Range("M61").Value = -1175
$ws.Range("N61").Value = -10516.2
$ws.Range("H82").Value = 2048.7273
$ws.Range("I82").Value = 1331.8334
$ws.Range("J82").Value = 2909
$ws.Range("K82").Value = 1331.8334
$ws.Range("L82").Value = 2909
$ws.Range("M82").Value = -970.8334
$ws.Range("N82").Value = -3631
$ws.Range("H85").Value = 2048.7273
$ws.Range("I85").Value = 1331.8334
$ws.Range("J85").Value = 2909
$ws.Range("K85").Value = 1331.8334
$ws.Range("L85").Value = 2909
$ws.Range("M85").Value = -83.83339999999998
$ws.Range("N85").Value = -5405
$ws.Range("H100").Value = 4743.0835
$ws.Range("I100").Value = 3713.4285
$ws.Range("J100").Value = 6184.6
$ws.Range("K100").Value = 3713.4285
$ws.Range("L100").Value = 6184.6
$ws.Range("M100").Value = -3172.4285
$ws.Range("N100").Value = -7266.6
$ws.Range("H108").Value = 67799.336
$ws.Range("J108").Value = 67799.336
$ws.Range("L108").Value = 67799.336
$ws.Range("N108").Value = -75479.336
$ws.Range("H113").Value = 3803.4443
$ws.Range("I113").Value = 1377
$ws.Range("J113").Value = 10112.2
$ws.Range("K113").Value = 1377
$ws.Range("L113").Value = 10112.2
$ws.Range("M113").Value = 793
$ws.Range("N113").Value = -14452.2
$ws.Range("H115").Value = 50000
$ws.Range("J115").Value = 50000
$ws.Range("L115").Value = 50000
$ws.Range("N115").Value = -52350
$ws.Range("H126").Value = 46090.75
$ws.Range("I126").Value = 62617
$ws.Range("J126").Value = 5955.5713
$ws.Range("K126").Value = 187851
$ws.Range("L126").Value = 17866.7139
$ws.Range("M126").Value = -185381
$ws.Range("N126").Value = -22806.7139
$ws.Range("H132").Value = 14524.053
$ws.Range("I132").Value = 13904.462
$ws.Range("J132").Value = 15866.5
$ws.Range("K132").Value = 41713.386
$ws.Range("L132").Value = 47599.5
$ws.Range("M132").Value = -39183.386
$ws.Range("N132").Value = -52659.5
$ws.Range("H136").Value = 8293
$ws.Range("I136").Value = 7374.778
$ws.Range("K136").Value = 22124.334
$ws.Range("M136").Value = -19574.334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1179.5
$ws.Range("I107").Value = 1274.2
$ws.Range("J107").Value = 1021.6667
$ws.Range("K107").Value = 3822.6
$ws.Range("L107").Value = 3065.0001
$ws.Range("M107").Value = -1902.6
$ws.Range("N107").Value = -6905.0001
$ws.Range("H113").Value = 523.4545000000001
$ws.Range("I113").Value = 577.375
$ws.Range("J113").Value = 379.66666
$ws.Range("K113").Value = 1732.125
$ws.Range("L113").Value = 1138.99998
$ws.Range("M113").Value = 437.875
$ws.Range("N113").Value = -5478.999980000001
$ws.Range("H126").Value = 4452.0605
$ws.Range("I126").Value = 4480.0347
$ws.Range("K126").Value = 13440.1041
$ws.Range("M126").Value = -10970.1041
$ws.Range("H132").Value = 4189.037
$ws.Range("I132").Value = 3624.9583
$ws.Range("J132").Value = 8701.666999999999
$ws.Range("K132").Value = 10874.8749
$ws.Range("L132").Value = 26105.001
$ws.Range("M132").Value = -8344.874899999999
$ws.Range("N132").Value = -31165.001
$ws.Range("H136").Value = 5728.193
$ws.Range("I136").Value = 3075.5898
$ws.Range("K136").Value = 9226.769400000001
$ws.Range("M136").Value = -6676.769400000001
